# fix for auth page
# Appends one new "submission time" measurement row to each of the four
# tracking sheets (mirrors new data collected on 10.13.2022).

$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Submit orders" -> append row 83 ---------------------------
$ws1 = $wb.Worksheets.Item("Submit orders")
$ws1.Range("A83:E83").Style = "Normal"
$ws1.Range("A83").Value = "10.13.2022 18:52 (Kyiv+Israel) 15:52 (UTC) 00:52 (Japan) 21:22 (India)"
$ws1.Range("B83").Value = "***"
$ws1.Range("C83").Value = "***"
$ws1.Range("D83").Value = 1.433
$ws1.Range("E83").Value = -0.2950000000000002

# --- Sheet 2: "Submit internet survey" -> append row 78 ------------------
$ws2 = $wb.Worksheets.Item("Submit internet survey")
$ws2.Range("A78:E78").Style = "Normal"
$ws2.Range("A78").Value = "10.13.2022 18:54 (Kyiv+Israel) 15:54 (UTC) 00:54 (Japan) 21:24 (India)"
$ws2.Range("B78").Value = "***"
$ws2.Range("C78").Value = "***"
$ws2.Range("D78").Value = 0.593
$ws2.Range("E78").Value = 0.172

# --- Sheet 3: "Submit a phone survey" -> append row 73 --------------------
$ws3 = $wb.Worksheets.Item("Submit a phone survey")
$ws3.Range("A73:E73").Style = "Normal"
$ws3.Range("A73").Value = "10.13.2022 18:57 (Kyiv+Israel) 15:57 (UTC) 00:57 (Japan) 21:27 (India)"
$ws3.Range("B73").Value = "***"
$ws3.Range("C73").Value = "***"
$ws3.Range("D73").Value = 1.917
$ws3.Range("E73").Value = -0.333

# --- Sheet 4: "Checkertificate" -> append rows 88 and 89 ------------------
$ws4 = $wb.Worksheets.Item("Checkertificate")

$ws4.Range("A88:E88").Style = "Normal"
$ws4.Range("A88").Value = "10.13.2022 18:59 (Kyiv+Israel) 15:59 (UTC) 00:59 (Japan) 21:29 (India)"
$ws4.Range("B88").Value = "***"
$ws4.Range("C88").Value = "***"
$ws4.Range("D88").Value = 0.996
$ws4.Range("E88").Value = -0.07399999999999995

$ws4.Range("A89:E89").Style = "Normal"
$ws4.Range("A89").Value = "10.13.2022 19:23 (Kyiv+Israel) 16:23 (UTC) 01:23 (Japan) 21:53 (India)"
$ws4.Range("B89").Value = 0.735
$ws4.Range("C89").Value = -0.06999999999999995
$ws4.Range("D89").Value = "***"
$ws4.Range("E89").Value = "***"
